$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.155946254730225
$ws.Range("B1").Value = 2.338923454284668
$ws.Range("C1").Value = 4.909433364868164
$ws.Range("D1").Value = 2.46947455406189
$ws.Range("E1").Value = 1.241406679153442
